$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 5 (shifts existing rows 5-8 down to 6-9)
$ws.Rows.Item(5).Insert()

# Row 4: rename label and change value (AG-Anteil)
$ws.Range("A4").Value = "Zusatzbeitrag Krankenversicherung AG-Anteil in Prozent"
$ws.Range("B4").Value = 0.99

# Row 5 (new row): AN-Anteil
$ws.Range("A5").Value = "Zusatzbeitrag Krankenversicherung AN-Anteil in Prozent"
$ws.Range("B5").Value = 0.99

# Row 6: Umlage U1 in Prozent value changed from 1.6 to 2.3
$ws.Range("B6").Value = 2.3

# Row 9: Eintragungsdatum value changed from 15.12.2023 to 01.01.2024
$ws.Range("B9").Value = "01.01.2024"

# Update selection to C9 as seen in the diff
$ws.Range("C9").Select()
